$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.118126630783081
$ws.Range("B1").Value = 3.592410326004028
$ws.Range("C1").Value = 3.270452499389648
$ws.Range("D1").Value = 3.662342309951782
$ws.Range("E1").Value = 1.067794442176819
